$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 8932
$ws1.Range("F13").Value = 1017
$ws1.Range("F17").Value = 241
$ws1.Range("F18").Value = 299
$ws1.Range("F19").Value = 72
$ws1.Range("F21").Value = 1105

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 8932
$ws4.Range("F15").Value = 1017
$ws4.Range("F19").Value = 241
$ws4.Range("F20").Value = 299
$ws4.Range("F21").Value = 72
$ws4.Range("F23").Value = 1105
